$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 74 updates
$ws.Range("B74").Value = 19402
$ws.Range("C74").Value = 17967
$ws.Range("D74").Value = 5803
$ws.Range("N74").Value = 2009
$ws.Range("O74").Value = 400
$ws.Range("R74").Value = 79
$ws.Range("S74").Value = 58
$ws.Range("V74").Value = 9543
$ws.Range("W74").Value = 2468
$ws.Range("X74").Value = 1042
$ws.Range("Y74").Value = 652
$ws.Range("AA74").Value = 235
$ws.Range("AB74").Value = 46
$ws.Range("AC74").Value = 95
$ws.Range("AD74").Value = 7075
$ws.Range("AF74").Value = 190
$ws.Range("AI74").Value = 711
$ws.Range("AK74").Value = 94

# Row 75 updates
$ws.Range("B75").Value = 20672
$ws.Range("C75").Value = 19060
$ws.Range("D75").Value = 6106
$ws.Range("F75").Value = 825
$ws.Range("H75").Value = 476
$ws.Range("I75").Value = 241
$ws.Range("J75").Value = 280
$ws.Range("K75").Value = 1365
$ws.Range("L75").Value = 620
$ws.Range("N75").Value = 2185
$ws.Range("O75").Value = 465
$ws.Range("P75").Value = 409
$ws.Range("U75").Value = 292
$ws.Range("V75").Value = 10279
$ws.Range("W75").Value = 2668
$ws.Range("X75").Value = 1058
$ws.Range("Y75").Value = 591
$ws.Range("AA75").Value = 400
$ws.Range("AC75").Value = 129
$ws.Range("AD75").Value = 7611
$ws.Range("AF75").Value = 258
$ws.Range("AG75").Value = 998
$ws.Range("AI75").Value = 645
$ws.Range("AJ75").Value = 100
$ws.Range("AO75").Value = 220
$ws.Range("AP75").Value = 4287
$ws.Range("AQ75").Value = 768
$ws.Range("AR75").Value = 108
$ws.Range("AT75").Value = 300
$ws.Range("AY75").Value = 1011
$ws.Range("AZ75").Value = 268
$ws.Range("BA75").Value = 277
